$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8429999947547913
$ws.Range("C2").Value = 0.830750048160553
$ws.Range("D2").Value = 0.8237500190734863
$ws.Range("E2").Value = 0.8260000348091125
$ws.Range("F2").Value = 0.8195000290870667

# Row 3
$ws.Range("B3").Value = 0.8882500529289246
$ws.Range("C3").Value = 0.8852500319480896
$ws.Range("D3").Value = 0.8795000314712524
$ws.Range("E3").Value = 0.8757500648498535
$ws.Range("F3").Value = 0.8782500028610229

# Row 4 - F4 is removed entirely
$ws.Range("B4").Value = 0.8486669063568115
$ws.Range("C4").Value = 0.8415879607200623
$ws.Range("D4").Value = 0.831055760383606
$ws.Range("E4").Value = 0.7680379152297974
$ws.Range("F4").ClearContents()

# Row 5 - F5 is removed entirely
$ws.Range("B5").Value = 0.8899509906768799
$ws.Range("C5").Value = 0.8891602754592896
$ws.Range("D5").Value = 0.881010890007019
$ws.Range("E5").Value = 0.853326141834259
$ws.Range("F5").ClearContents()

# Row 6
$ws.Range("B6").Value = 0.5597222447395325
$ws.Range("C6").Value = 0.7240484356880188
$ws.Range("D6").Value = 0.8057762384414673
$ws.Range("E6").Value = 0.8247553110122681
$ws.Range("F6").Value = 0.8173292279243469

# Row 7
$ws.Range("B7").Value = 0.8159722089767456
$ws.Range("C7").Value = 0.8330006003379822
$ws.Range("D7").Value = 0.8664716482162476
$ws.Range("E7").Value = 0.8791993856430054
$ws.Range("F7").Value = 0.8810346722602844

# Row 8
$ws.Range("B8").Value = 0.8537499904632568
$ws.Range("C8").Value = 0.5720000267028809
$ws.Range("D8").Value = 0.4925000071525574
$ws.Range("E8").Value = 0.8295000791549683
$ws.Range("F8").Value = 0.8715000152587891

# Row 9
$ws.Range("B9").Value = 0.01750000193715096
$ws.Range("C9").Value = 0.101000003516674
$ws.Range("D9").Value = 0.4147500395774841
$ws.Range("E9").Value = 0.8287500143051147
$ws.Range("F9").Value = 0.8715000152587891

# Row 10
$ws.Range("B10").Value = 0.7559899091720581
$ws.Range("C10").Value = 0.7352020740509033
$ws.Range("D10").Value = 0.7666065096855164
$ws.Range("E10").Value = 0.797055721282959
$ws.Range("F10").Value = 0.7952351570129395

# Row 11 - E11 is removed entirely
$ws.Range("B11").Value = 0.8035588264465332
$ws.Range("C11").Value = 0.8076726198196411
$ws.Range("D11").Value = 0.8219923377037048
$ws.Range("E11").ClearContents()
$ws.Range("F11").Value = 0.8484794497489929

# Row 12 - B12 is a newly added cell
$ws.Range("B12").Value = 0.7559899091720581
$ws.Range("C12").Value = 0.7352020740509033
$ws.Range("D12").Value = 0.7666065096855164
$ws.Range("E12").Value = 0.797055721282959
$ws.Range("F12").Value = 0.8115577697753906
